$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (tckn, isim-soyisim, Rol) replacing the old sample rows.
# The last flag marks whether column C keeps its highlighted ("style 2")
# formatting for that row, matching the source data.
$data = @(
    @(12345678929, "FULYA İNCİ", "Takım üyesi,Yazılımcı", $true),
    @(12345678930, "ÖMER FARUK BORAN", "Yazılımcı", $true),
    @(12345678931, "MELİKE ERTAN", "Yardımcı", $false),
    @(12345678932, "MUHAMMED ALİ KÖSEN", "Araştırmacı,Yardımcı", $true),
    @(12345678933, "MUHAMMED GÖNEN", "Yazılımcı,Takım üyesi", $true),
    @(12345678934, "FATMANUR ÖZDEMİR", "Yardımcı", $false)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]

    # The refreshed rows lose the explicit "s=1" style that the old sample
    # rows carried (Excel re-applied the plain default style when the
    # values were replaced). Column C keeps its own highlighted style for
    # most rows, but a couple of rows drop back to the default too.
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Style = "Normal"
    if (-not $entry[3]) {
        $ws.Cells.Item($row, 3).Style = "Normal"
    }

    $row = $row + 1
}
